# Task: Задание 3.4.2, 3.4.3, 3.5.1
# Updates the "years" statistics sheet (sheet 1) with a new data range (rows 2-21,
# years 2003-2022) and the "cities" statistics sheet (sheet 2) with refreshed
# salary-level and vacancy-share rankings, plus adjusted column widths.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# Sheet 1 ("Статистика по годам")
# ---------------------------------------------------------------------------

# Extend the data block (previously rows 2-17) down to row 21 by copying the
# formatting (style/borders/alignment) of the last existing data row.
$ws1.Range("A17:E17").Copy($ws1.Range("A18:E21"))

$years1 = @(2003,2004,2005,2006,2007,2008,2009,2010,2011,2012,2013,2014,2015,2016,2017,2018,2019,2020,2021,2022)
$colB   = @(1366,1488,13331,1522,5604,27478,37548,40958,42359,44540,46218,48482,50654,58261,61724,65563,78212,90537,105356,124935)
$colC   = @($null,$null,$null,$null,$null,$null,$null,$null,40000,31250,25000,27500,28625,26800,32681,33250,34708,47500,50000,$null)
$colD   = @(1983,7833,16022,33321,53562,75070,52889,93494,142458,173897,234019,259571,284763,332460,391464,517670,535956,489472,287915,91142)
$colE   = @(0,0,0,0,0,0,0,0,2,5,7,1,5,14,13,18,17,9,3,1)

for ($i = 0; $i -lt $years1.Count; $i++) {
    $r = $i + 2
    $ws1.Cells.Item($r,1).Value = $years1[$i]
    $ws1.Cells.Item($r,2).Value = $colB[$i]
    if ($null -eq $colC[$i]) {
        $ws1.Cells.Item($r,3).ClearContents()
    } else {
        $ws1.Cells.Item($r,3).Value = $colC[$i]
    }
    $ws1.Cells.Item($r,4).Value = $colD[$i]
    $ws1.Cells.Item($r,5).Value = $colE[$i]
}

# ---------------------------------------------------------------------------
# Sheet 2 ("Статистика по городам")
# ---------------------------------------------------------------------------

$colA2 = @("Алматы","Москва","Санкт-Петербург","Новосибирск","Екатеринбург","Краснодар","Казань","Самара","Нижний Новгород","Пермь")
$colB2 = @(179001,70218,61230,57592,54972,50049,47800,46337,44775,44542)
$colD2 = @("Москва","Санкт-Петербург","Минск","Киев","Новосибирск","Нижний Новгород","Екатеринбург","Алматы","Воронеж","Казань")
$colE2 = @(0.4917,0.1493,0.0598,0.0474,0.0348,0.0316,0.029,0.0279,0.0274,0.0273)

for ($i = 0; $i -lt $colA2.Count; $i++) {
    $r = $i + 2
    $ws2.Cells.Item($r,1).Value = $colA2[$i]
    $ws2.Cells.Item($r,2).Value = $colB2[$i]
    $ws2.Cells.Item($r,4).Value = $colD2[$i]
    $ws2.Cells.Item($r,5).Value = $colE2[$i]
}

# Column widths on sheet 2 grow by one character each (raw OOXML width =
# ColumnWidth + 5/6, so subtract 5/6 from the desired raw widths).
$offset = 5/6
$ws2.Columns.Item(1).ColumnWidth = 17 - $offset
$ws2.Columns.Item(2).ColumnWidth = 17 - $offset
$ws2.Columns.Item(3).ColumnWidth = 4 - $offset
$ws2.Columns.Item(4).ColumnWidth = 17 - $offset
$ws2.Columns.Item(5).ColumnWidth = 15 - $offset

Write-Host "Edit applied"
